$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.868833184242249
$ws.Range("B1").Value = 2.024146556854248
$ws.Range("C1").Value = 2.010165691375732
$ws.Range("D1").Value = 2.404569149017334
$ws.Range("E1").Value = 2.639673709869385
